$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): add Wins, Losses, Ties in AD1:AF1.
# Copy the formatting of an existing bold/bordered/centered header cell (AC1)
# onto the new header cells, then overwrite their text.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-37: season record (Wins=116, Losses=46, Ties=0) repeated for each player
for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 30).Value = 116   # AD
    $ws.Cells.Item($r, 31).Value = 46    # AE
    $ws.Cells.Item($r, 32).Value = 0     # AF
}
